# Create Engineering Item Master.xlsx - QA automation data maintenance
# (Added test cases for Recurring billing / maintenance of other test cases)

$wb = $excel.ActiveWorkbook

# --- "Create Engg Item" sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Create Engg Item")

# Inventory Division column (C2) now points at a different test site value.
$ws1.Range("C2").Value = "Colorado (100)"

# Column C widened (Excel "best fit") to accommodate the new, longer value.
$ws1.Columns.Item(3).ColumnWidth = 12.39

# Cursor was left on I19 when the sheet was last touched/saved.
$ws1.Range("I19").Select()

# --- "Routing Master" sheet ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Routing Master")

# New recurring-billing style test item pushed through Item Number / Item
# Description / Id, recording a fresh set of Provar/Salesforce test ids.
$ws2.Range("B2").Value = "Pro-PEItem-MFMOK"
$ws2.Range("C2").Value = "Pro-PEItem-H58GD"
$ws2.Range("D2").Value = "a345f000000uGVNAA2"

# Columns re-sized ("best fit") because of the new, longer values.
$ws2.Columns.Item(1).ColumnWidth = 10.72
$ws2.Columns.Item(2).ColumnWidth = 18.34
$ws2.Columns.Item(3).ColumnWidth = 17.31
$ws2.Columns.Item(4).ColumnWidth = 20.83

# "Routing Master" was the active tab when the workbook was saved.
$ws2.Activate()
